# Append a new data row (row 3) to the "Artfynd" sheet, mirroring the
# structure of the existing rows (row 1 = headers, row 2 = first record).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- numeric columns -------------------------------------------------
$ws.Range("A3").Value = 111977935
$ws.Range("B3").Value = 88032
$ws.Range("E3").Value = 6276
$ws.Range("Q3").Value = 538762.8579659602
$ws.Range("R3").Value = 6718247.433583082
$ws.Range("S3").Value = 5

# --- plain text columns ------------------------------------------------
$ws.Range("C3").Value = "Ovaliderad"
$ws.Range("D3").Value = "VU"
$ws.Range("F3").Value = "Goliatmusseron"
$ws.Range("G3").Value = "Tricholoma matsutake"
$ws.Range("H3").Value = "(S.Ito & S.Imai) Singer"
$ws.Range("J3").Value = "fruktkroppar"
$ws.Range("P3").Value = "Hälsingbergsskogen, Dlr"
$ws.Range("T3").Value = "Dalarna"
$ws.Range("U3").Value = "Falun"
$ws.Range("V3").Value = "Dalarna"
$ws.Range("W3").Value = "Stora Kopparberg"
$ws.Range("AC3").Value = "Fruktkroppen övermogen. Lavrik stenig tallskog."
$ws.Range("AW3").Value = "Uno Skog"
$ws.Range("AX3").Value = "Uno Skog"

# --- text columns that look numeric/date-like: force text with a
# leading apostrophe so Excel does not auto-convert them -------------
$ws.Range("I3").Value = "'1"
$ws.Range("Y3").Value = "'2023-09-09"
$ws.Range("Z3").Value = "'10:53"
$ws.Range("AA3").Value = "'2023-09-09"
$ws.Range("AB3").Value = "'10:53"

# --- boolean columns ---------------------------------------------------
$ws.Range("AD3").Value = $false
$ws.Range("AE3").Value = $false
$ws.Range("AG3").Value = $false

# --- empty (but present) text cells ------------------------------------
$ws.Range("K3").Value = "'"
$ws.Range("AT3").Value = "'"
$ws.Range("AY3").Value = "'"
